$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (dates are Excel serial numbers, continuing the daily series)
# Columns: A = date, B = nuovi pos., C = somma mobile 7gg., D = somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44403, 0, 0, 0),
    @(44404, 0, 0, 0),
    @(44405, 0, 0, 0),
    @(44406, 0, 0, 0),
    @(44407, 0, 0, 0),
    @(44408, 0, 0, 0),
    @(44409, 2, 2, 74.93443237167479),
    @(44410, 0, 2, 74.93443237167479),
    @(44411, 0, 2, 74.93443237167479),
    @(44412, 0, 2, 74.93443237167479),
    @(44413, 0, 2, 74.93443237167479),
    @(44414, 0, 2, 74.93443237167479),
    @(44415, 0, 2, 74.93443237167479),
    @(44416, 0, 0, 0),
    @(44417, 0, 0, 0)
)

$startRow = 329
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Column A: date value, styled like the rest of column A (copy format from row above)
    $ws.Range("A$($r - 1)").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Cells.Item($r, 1).Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = 0
